$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Row 20 ("cat", first block, rows 1-22) ---
$row20 = @(
    0.344247787610619,
    0.299115044247787,
    0.308849557522123,
    0.266371681415929,
    0.207079646017699,
    0.179646017699115,
    0.136283185840707,
    0.131858407079646,
    0.0964601769911504,
    0.0336283185840707,
    0.0008849557522123
)

# --- Row 21 ("   +all", first block, rows 1-22) ---
$row21 = @(
    0.631858407079646,
    0.569911504424778,
    0.504424778761062,
    0.497345132743362,
    0.367256637168141,
    0.346017699115044,
    0.286725663716814,
    0.227433628318584,
    0.181415929203539,
    0.117699115044247,
    0.0628318584070796
)

# --- Row 43 ("cat", second block, rows 23-45) ---
$row43 = @(
    0.451914354770822,
    0.400939273435206,
    0.394813672236949,
    0.349617502553614,
    0.278921125800195,
    0.236864714172954,
    0.179863706412049,
    0.155327595672601,
    0.134506537493539,
    0.046163942497512,
    0.0081632858704344
)

# --- Row 44 ("   +all", second block, rows 23-45) ---
$row44 = @(
    0.74347147090489,
    0.67664885963833,
    0.604584714475048,
    0.566152904867849,
    0.449170940140829,
    0.427331699735831,
    0.353471855374924,
    0.299201748851832,
    0.234046338367176,
    0.161633927837739,
    0.100092752284885
)

for ($i = 0; $i -lt $row20.Length; $i++) {
    $col = 2 + $i   # Column B = 2
    $ws.Cells.Item(20, $col).Value = $row20[$i]
    $ws.Cells.Item(21, $col).Value = $row21[$i]
    $ws.Cells.Item(43, $col).Value = $row43[$i]
    $ws.Cells.Item(44, $col).Value = $row44[$i]
}

$excel.Calculate()

# --- Update the active selection on Sheet2 ---
$ws.Activate()
[void]$ws.Range("T30").Select()
